{"js": "// The document's built-in \"Footnote Text\" paragraph style (styleId\n// \"FootnoteText\") is unused in the body/footnotes content and is being\n// removed from the style sheet entirely. Locate it by its display name\n// and delete it.\nconst styles = context.document.getStyles();\nconst footnoteTextStyle = styles.getByNameOrNullObject(\"Footnote Text\");\nawait context.sync();\n\nif (!footnoteTextStyle.isNullObject) {\n  footnoteTextStyle.delete();\n  await context.sync();\n}\n", "ps1": "# The document's built-in \"Footnote Text\" paragraph style (styleId\n# \"FootnoteText\") is unused in the body/footnotes content and is being\n# removed from the style sheet entirely. Locate it by its display name\n# and delete it.\n$d = $word.ActiveDocument\n\ntry {\n    $footnoteTextStyle = $d.Styles(\"Footnote Text\")\n} catch {\n    $footnoteTextStyle = $null\n}\n\nif ($footnoteTextStyle -ne $null) {\n    $footnoteTextStyle.Delete()\n}\n"}
